# Insert a new "budget-type" column into the "Data-wide-value" sheet.
#
# Before: A=entity-name, B..F = year columns (2012..2016)
# After:  A=entity-name, B="budget-type" (all rows = "budget"), C..G = year columns (2012..2016)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Shift existing columns B:F to C:G by inserting a new, blank column at B.
$ws.Columns("B:B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "budget-type"

# Determine how many data rows exist (row 1 is the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Fill every data row's new column B with the constant "budget".
$ws.Range("B2:B" + $lastRow).Value = "budget"
